$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (data rows) and append a new row right after it
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Force column A to be treated as plain text so the date string is not
# auto-converted into a date serial number, then reset the style back to
# the default ("Normal") so no extra style index is left attached to the cell.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025-10-12"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 54.31000137329102
$ws.Cells.Item($newRow, 3).Value = 678.9500122070312
$ws.Cells.Item($newRow, 4).Value = 348.2999877929688
